$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.919.85'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.83%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.51'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.33%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '281.41'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.19%  '

$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5073'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3523'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.69%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.23'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.48%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06646'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.12'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8529'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.39%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07855'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.822.40'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +63.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.038'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.38'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.76%  '

$ws.Range('E17').Value = '  +0.44%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.04'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008137'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E20').Value = '  +0.20%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '25.980.42'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.67%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.762'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.35%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.11'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.120'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.57'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.43%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.152'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.88%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.678'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.88%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.95'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.56%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '108.54'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.93%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.315'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.213'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -11.84%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08792'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.92%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04795'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7407'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.76%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.125'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.80%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.854'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.67%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.003'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.102'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.438'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.87%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5388'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01854'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.38%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9857'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.42%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '112.84'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.67%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.237'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.213'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -11.97%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4705'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.68%  '

$ws.Range('E47').Value = '  +0.22%  '

$ws.Range('E48').Value = '  -9.46%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.240'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.76'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.63%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05911'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.66%  '
